$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number but must stay text
# (mirrors the source inlineStr cells), so force Text number format first.
$numericLookingCells = @("D4", "D5", "D6", "D11", "D13", "D15", "D19", "D20", "D21", "D22", "D23", "D24", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '67.118.05'
$ws.Range('E2').Value = '  -3.37%  '
$ws.Range('D3').Value = '3.538.48'
$ws.Range('E3').Value = '  -3.82%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '608.73'
$ws.Range('E5').Value = '  -5.38%  '
$ws.Range('D6').Value = '154.28'
$ws.Range('E6').Value = '  -2.88%  '
$ws.Range('D7').Value = '3.534.45'
$ws.Range('E7').Value = '  -3.85%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -2.71%  '
$ws.Range('E10').Value = '  -2.69%  '
$ws.Range('D11').Value = '6.88'
$ws.Range('E11').Value = '  -2.81%  '
$ws.Range('E12').Value = '  -3.47%  '
$ws.Range('D13').Value = '0.0000221'
$ws.Range('E13').Value = '  -4.46%  '
$ws.Range('D14').Value = '4.130.93'
$ws.Range('E14').Value = '  -3.83%  '
$ws.Range('D15').Value = '31.89'
$ws.Range('E15').Value = '  -1.95%  '
$ws.Range('D16').Value = '3.532.92'
$ws.Range('E16').Value = '  -3.98%  '
$ws.Range('D17').Value = '67.010.52'
$ws.Range('E17').Value = '  -3.43%  '
$ws.Range('E18').Value = '  +0.44%  '
$ws.Range('D19').Value = '6.32'
$ws.Range('E19').Value = '  -2.12%  '
$ws.Range('D20').Value = '15.42'
$ws.Range('E20').Value = '  -3.28%  '
$ws.Range('D21').Value = '447.81'
$ws.Range('E21').Value = '  -4.09%  '
$ws.Range('D22').Value = '9.28'
$ws.Range('E22').Value = '  -5.77%  '
$ws.Range('D23').Value = '0.633'
$ws.Range('E23').Value = '  -1.91%  '
$ws.Range('D24').Value = '78.31'
$ws.Range('E24').Value = '  -1.21%  '
$ws.Range('D25').Value = '3.675.42'
$ws.Range('E25').Value = '  -3.85%  '
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('D27').Value = '0.0000123'
$ws.Range('E27').Value = '  -1.57%  '
$ws.Range('D28').Value = '10.27'
$ws.Range('E28').Value = '  -4.85%  '
$ws.Range('D29').Value = '8.23'
$ws.Range('E29').Value = '  -8.57%  '
$ws.Range('D30').Value = '2.54'
$ws.Range('E30').Value = '  -3.30%  '
$ws.Range('D31').Value = '1.66'
$ws.Range('E31').Value = '  -1.87%  '
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  -0.20%  '
$ws.Range('D33').Value = '25.76'
$ws.Range('E33').Value = '  -4.04%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = '1.89'
$ws.Range('E34').Value = '  -5.23%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = '0.158'
$ws.Range('E35').Value = '  -4.08%  '
$ws.Range('D36').Value = '6.17'
$ws.Range('E36').Value = '  -3.73%  '
$ws.Range('D37').Value = '3.528.33'
$ws.Range('E37').Value = '  -3.82%  '
$ws.Range('D38').Value = '8.03'
$ws.Range('E38').Value = '  -4.33%  '
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('D40').Value = '0.998'
$ws.Range('E40').Value = '  -0.07%  '
$ws.Range('D41').Value = '176.02'
$ws.Range('E41').Value = '  -0.70%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '2.15'
$ws.Range('E42').Value = '  -3.67%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').Value = '5.58'
$ws.Range('E43').Value = '  -5.00%  '
$ws.Range('D44').Value = '0.0867'
$ws.Range('E44').Value = '  -3.09%  '
$ws.Range('D45').Value = '0.891'
$ws.Range('E45').Value = '  -3.59%  '
$ws.Range('D46').Value = '45.75'
$ws.Range('E46').Value = '  -2.24%  '
$ws.Range('D47').Value = '27.93'
$ws.Range('E47').Value = '  +2.32%  '
$ws.Range('D48').Value = '2.63'
$ws.Range('E48').Value = '  -2.49%  '
$ws.Range('D49').Value = '1.23'
$ws.Range('E49').Value = '  -1.09%  '
$ws.Range('D50').Value = '7.60'
$ws.Range('E50').Value = '  -2.71%  '
$ws.Range('E51').Value = '  -3.45%  '
